$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was updated
# from 45180 (2023-09-11) to 45181 (2023-09-12) for every data row
# (rows 2 through 236).
$ws.Range("C2:C236").Value = 45181
